# Slide 1 ("sldId 256"): the "Entrance / Time / Payment type" table
# (shape id 29, a:tbl) had its first column narrowed, which in turn
# shrinks the overall graphic frame width (height/position unchanged).
#
# EMU values taken from the target OOXML:
#   first column width: 1116000 EMU -> 743268 EMU   (delta -372732 EMU)
#   frame cx:           3348000 EMU -> 2975268 EMU  (same delta)
# PowerPoint's object model works in points (1 pt = 12700 EMU), and
# resizing Columns.Item(1).Width automatically adjusts the parent
# Shape.Width by the same delta, leaving the other columns untouched.

$EMU_PER_POINT = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 29) {
        $shp = $candidate
        break
    }
}

$table = $shp.Table
$firstCol = $table.Columns.Item(1)
$firstCol.Width = 743268 / $EMU_PER_POINT
